# Update "Pais" sheet with refreshed COVID-19 country statistics
# and re-ranked rows (countries that changed order by total cases).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp update
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 18:59"

# Row 4: refreshed stats
$ws.Range("B4").Value = 5213945
$ws.Range("C4").Value = 14501
$ws.Range("D4").Value = 2667704
$ws.Range("E4").Value = 2380459
$ws.Range("G4").Value = 165
$ws.Range("H4").Value = 165782

# Row 6: refreshed stats
$ws.Range("B6").Value = 2257572
$ws.Range("C6").Value = 43435
$ws.Range("D6").Value = 1571712
$ws.Range("E6").Value = 640662
$ws.Range("G6").Value = 732
$ws.Range("H6").Value = 45198

# Row 22: refreshed stats
$ws.Range("B22").Value = 217330
$ws.Range("C22").Value = 49
$ws.Range("E22").Value = 10170

# Row 59: refreshed stats
$ws.Range("B59").Value = 35712
$ws.Range("C59").Value = 498
$ws.Range("D59").Value = 24920
$ws.Range("E59").Value = 9480
$ws.Range("G59").Value = 10
$ws.Range("H59").Value = 1312

# Row 63: refreshed stats
$ws.Range("D63").Value = 18965
$ws.Range("E63").Value = 8651

# Row 66: refreshed stats
$ws.Range("B66").Value = 26768
$ws.Range("C66").Value = 56
$ws.Range("E66").Value = 1632

# Row 68: country changed to "Etiopia" with refreshed stats
$ws.Range("A68").Value = "Etiopia"
$ws.Range("B68").Value = 23591
$ws.Range("C68").Value = 773
$ws.Range("D68").Value = 10411
$ws.Range("E68").Value = 12760
$ws.Range("G68").Value = 13
$ws.Range("H68").Value = 420

# Row 69: country changed to "Nepal" with refreshed stats
$ws.Range("A69").Value = "Nepal"
$ws.Range("B69").Value = 23310
$ws.Range("C69").Value = 338
$ws.Range("D69").Value = 16493
$ws.Range("E69").Value = 6738
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 79

# Row 70: country changed to "Costa Rica" with refreshed stats
$ws.Range("A70").Value = "Costa Rica"
$ws.Range("B70").Value = 23286
$ws.Range("D70").Value = 7730
$ws.Range("E70").Value = 15321
$ws.Range("H70").Value = 235

# Row 74: refreshed stats
$ws.Range("B74").Value = 18454
$ws.Range("C74").Value = 101
$ws.Range("D74").Value = 12983
$ws.Range("E74").Value = 5081

# Row 84: refreshed stats
$ws.Range("B84").Value = 11942
$ws.Range("C84").Value = 103
$ws.Range("D84").Value = 8087
$ws.Range("E84").Value = 3327
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 528

# Row 96: refreshed stats
$ws.Range("B96").Value = 7216
$ws.Range("C96").Value = 11
$ws.Range("D96").Value = 6170
$ws.Range("E96").Value = 925
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 121

# Row 101: country changed to "Grecia" with refreshed stats
$ws.Range("A101").Value = "Grecia"
$ws.Range("B101").Value = 5749
$ws.Range("C101").Value = 126
$ws.Range("D101").Value = 3804
$ws.Range("E101").Value = 1732
$ws.Range("H101").Value = 213

# Row 102: country changed to "Croacia" with refreshed stats
$ws.Range("A102").Value = "Croacia"
$ws.Range("B102").Value = 5649
$ws.Range("C102").Value = 45
$ws.Range("D102").Value = 4906
$ws.Range("E102").Value = 585
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 158

# Row 121: country changed to "Sri Lanka" with refreshed stats
$ws.Range("A121").Value = "Sri Lanka"
$ws.Range("B121").Value = 2869
$ws.Range("C121").Value = 25
$ws.Range("D121").Value = 2593
$ws.Range("E121").Value = 265
$ws.Range("H121").Value = 11

# Row 122: country changed to "Cabo Verde" with refreshed stats
$ws.Range("A122").Value = "Cabo Verde"
$ws.Range("B122").Value = 2858
$ws.Range("D122").Value = 2086
$ws.Range("E122").Value = 740
$ws.Range("H122").Value = 32

# Row 126: country changed to "Mozambique" with refreshed stats
$ws.Range("A126").Value = "Mozambique"
$ws.Range("B126").Value = 2411
$ws.Range("C126").Value = 142
$ws.Range("D126").Value = 860
$ws.Range("E126").Value = 1535
$ws.Range("H126").Value = 16

# Row 127: country changed to "Surinam" with refreshed stats
$ws.Range("A127").Value = "Surinam"
$ws.Range("B127").Value = 2391
$ws.Range("D127").Value = 1635
$ws.Range("E127").Value = 727
$ws.Range("H127").Value = 29

# Row 135: refreshed stats
$ws.Range("B135").Value = 1917
$ws.Range("C135").Value = 1
$ws.Range("D135").Value = 1447
$ws.Range("E135").Value = 401
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 69

# Row 143: refreshed stats
$ws.Range("B143").Value = 1268
$ws.Range("C143").Value = 16
$ws.Range("E143").Value = 70

# Row 158: refreshed stats
$ws.Range("E158").Value = 434
$ws.Range("G158").Value = 3
$ws.Range("H158").Value = 14

# Row 202: country re-ranked (name only, stats follow below)
$ws.Range("A202").Value = "Timor Oriental"

# Row 203: country re-ranked (name only, stats follow below)
$ws.Range("A203").Value = "Santa Lucia"

# Row 209: refreshed stats
$ws.Range("D209").Value = 17
$ws.Range("E209").Value = 0

# Row 213: country changed to "Montserrat" with refreshed stats
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214: country changed to "Islas Malvinas" with refreshed stats
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
